$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.869.27'
$ws.Range("E2").Value = '  +2.64%  '
$ws.Range("D3").Value = '1.668.49'
$ws.Range("E3").Value = '  -0.47%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.85'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  -0.37%  '
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '23.57'
$ws.Range("E8").Value = '  +3.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.261'
$ws.Range("E9").Value = '  -0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0879'
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").Value = '1.904.58'
$ws.Range("E12").Value = '  -0.43%  '
$ws.Range("D13").Value = '1.668.82'
$ws.Range("E13").Value = '  -0.57%  '
$ws.Range("E14").Value = '  -1.44%  '
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.10'
$ws.Range("E16").Value = '  -0.63%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '251.67'
$ws.Range("E17").Value = '  +6.96%  '
$ws.Range("D18").Value = '27.833.80'
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("E19").Value = '  -1.16%  '
$ws.Range("E20").Value = '  -3.34%  '
$ws.Range("E21").Value = '  -0.13%  '
$ws.Range("E22").Value = '  -1.28%  '
$ws.Range("E23").Value = '  -2.24%  '
$ws.Range("E24").Value = '  -1.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.05'
$ws.Range("E25").Value = '  -0.88%  '
$ws.Range("E26").Value = '  -3.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.32'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  -0.43%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("E30").Value = '  +5.93%  '
$ws.Range("E31").Value = '  +0.29%  '
$ws.Range("E32").Value = '  -0.82%  '
$ws.Range("E33").Value = '  -2.98%  '
$ws.Range("D34").Value = '1.430.35'
$ws.Range("E34").Value = '  -7.43%  '
$ws.Range("E35").Value = '  -5.87%  '
$ws.Range("E36").Value = '  -0.01%  '
$ws.Range("E37").Value = '  -1.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.583'
$ws.Range("E38").Value = '  -3.99%  '
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("E40").Value = '  -2.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '69.84'
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.23'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.39'
$ws.Range("E44").Value = '  -6.66%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.812.34'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("E46").Value = '  +0.87%  '
$ws.Range("E47").Value = '  +4.55%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.15'
$ws.Range("E48").Value = '  -0.64%  '
$ws.Range("B49").Value = 'BabyDogeCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D49").Value = '0.0₆0111'
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.102'
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.86'
$ws.Range("E51").Value = '  -4.16%  '
